$d = $word.ActiveDocument

# The bullet list originally read:
#   2) "Calcular tensiones y direcciones principales del punto donde se
#       maximizan las tensiones en los mapas"
#   3) "A" + "nalíticamente, calcular las tensiones y direcciones
#       principales del punto a 1cm a la derecha del agujero"
#   4) "Cuanta fuerza se necesita para realizar la carga máxima
#       encontrada teóricamente" + ", para poder comparar y analizar
#       los resultados"
#
# The new instructions drop bullet 2 entirely (its content is replaced by
# what used to be bullet 3, which in turn is replaced by what used to be
# bullet 4), and append two brand-new bullets after the old bullet 4.

# 1) Remove the first bullet ("Calcular tensiones ..."); this shifts the
#    following paragraphs up by one slot, so the old bullets 3 and 4 keep
#    their original runs completely untouched.
$d.Paragraphs.Item(2).Range.Delete()

# 2) After the (now shifted up) "Cuanta fuerza ..." bullet, add a new
#    bullet with the new text about obtaining the tensor state at 1cm.
$pCuantaFuerza = $d.Paragraphs.Item(3)
$pCuantaFuerza.Range.InsertParagraphAfter()
$pEstadoTensorial = $d.Paragraphs.Item(4)
$pEstadoTensorial.Range.Text = "Para el estado tensorial a 1cm obtener las tensiones que se están aplicando a partir de la matriz que forma los mapas, así podemos saber s_xx, s_yy y t_xy. Y después podemos hacer el círculo de Mohr."

# 3) Add one more bullet after that, about finding the position of the
#    maximum value in the matrix.
$pEstadoTensorial2 = $d.Paragraphs.Item(4)
$pEstadoTensorial2.Range.InsertParagraphAfter()
$pValorMaximo = $d.Paragraphs.Item(5)
$pValorMaximo.Range.Text = "Igual para el valor máximo, encontrar una función que nos diga en qué posición se encuentra el valor máximo de la matriz y luego hacer el círculo de Mohr."
